$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 34: A34 = date serial 46003 (formatted like A33), B34 = 10
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A34").Value = 46003
$ws.Range("B34").Value = 10

# Update the selected/active cell to match the new last row
$ws.Range("A34:B34").Select()
